$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 00:45"
$ws.Range("B9").Value = 13816
$ws.Range("C9").Value = 4557
$ws.Range("E9").Value = 13501
$ws.Range("A24").Value = "Australia"
$ws.Range("B24").Value = 756
$ws.Range("C24").Value = 160
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = 704
$ws.Range("F24").Value = 1
$ws.Range("H24").Value = 6
$ws.Range("A25").Value = "Crucero"
$ws.Range("B25").Value = 712
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 527
$ws.Range("E25").Value = 178
$ws.Range("F25").Value = 14
$ws.Range("H25").Value = 7
$ws.Range("A102").Value = "Liechtenstein"
$ws.Range("C102").Value = 0
$ws.Range("A103").Value = "Nueva Zelanda"
$ws.Range("C103").Value = 8
$ws.Range("A111").Value = "Bolivia"
$ws.Range("C111").Value = 3
$ws.Range("A112").Value = "Guayana Francesa"
$ws.Range("C112").Value = 0
$ws.Range("A115").Value = "Montenegro"
$ws.Range("C115").Value = 5
$ws.Range("A116").Value = "Maldivas"
$ws.Range("A117").Value = "Camerun"
$ws.Range("C117").Value = 0
$ws.Range("A121").Value = "Ghana"
$ws.Range("C121").Value = 4
$ws.Range("A122").Value = "Ruanda"
$ws.Range("C122").Value = 0
$ws.Range("A128").Value = "Guatemala"
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 0
$ws.Range("H128").Value = 1
$ws.Range("A129").Value = "Costa de Marfil"
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 1
$ws.Range("H129").Value = 0
$ws.Range("A131").Value = "Kenia"
$ws.Range("C131").Value = 0
$ws.Range("A132").Value = "Mauricio"
$ws.Range("C132").Value = 4
$ws.Range("A134").Value = "Guinea Ecuatorial"
$ws.Range("C134").Value = 2
$ws.Range("A135").Value = "Seychelles"
$ws.Range("C135").Value = 0
$ws.Range("A136").Value = "Polinesia Francesa"
$ws.Range("C136").Value = 1
$ws.Range("A137").Value = "Tanzania"
$ws.Range("C137").Value = 3
$ws.Range("A138").Value = "Mongolia"
$ws.Range("A143").Value = "Gabon"
$ws.Range("A144").Value = "San Bartolome"
$ws.Range("A145").Value = "Congo"
$ws.Range("C145").Value = 2
$ws.Range("A146").Value = "Bahamas"
$ws.Range("C146").Value = 2
$ws.Range("A147").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("C147").Value = 1
$ws.Range("A148").Value = "Namibia"
$ws.Range("C148").Value = 1
$ws.Range("A149").Value = "San Martin (Parte Francesa)"
$ws.Range("C149").Value = 0
$ws.Range("A150").Value = "Kirguistan"
$ws.Range("C150").Value = 0
$ws.Range("A151").Value = "Curazao"
$ws.Range("C151").Value = 0
$ws.Range("G151").Value = 1
$ws.Range("A152").Value = "Islas Caimanes"
$ws.Range("C152").Value = 2
$ws.Range("G152").Value = 0
$ws.Range("A153").Value = "Groenlandia"
$ws.Range("A154").Value = "Bermudas"
$ws.Range("A156").Value = "Nueva Caledonia"
$ws.Range("A157").Value = "Santa Lucia"
$ws.Range("A159").Value = "Benin"
$ws.Range("A160").Value = "Zambia"
$ws.Range("A162").Value = "Somalia"
$ws.Range("A163").Value = "San Martin (Parte Holandesa)"
$ws.Range("A164").Value = "Republica de Yibuti"
$ws.Range("A165").Value = "Butan"
$ws.Range("A166").Value = "Niger"
$ws.Range("C166").Value = 1
$ws.Range("A167").Value = "Santa Sede"
$ws.Range("C167").Value = 0
$ws.Range("A168").Value = "Isla de Man"
$ws.Range("C168").Value = 1
$ws.Range("A169").Value = "Fiyi"
$ws.Range("C169").Value = 1
$ws.Range("A170").Value = "San Vicente y las Granadinas"
$ws.Range("A171").Value = "Suazilandia"
$ws.Range("A172").Value = "Gambia"
$ws.Range("A173").Value = "Montserrat"
$ws.Range("C173").Value = 0
$ws.Range("A174").Value = "Surinam"
$ws.Range("C174").Value = 0
$ws.Range("A175").Value = "Guinea"
$ws.Range("A176").Value = "Antigua y Barbuda"
$ws.Range("A177").Value = "El Salvador"
$ws.Range("C177").Value = 1
$ws.Range("A178").Value = "Nicaragua"
$ws.Range("C178").Value = 1
$ws.Range("A179").Value = "Togo"
$ws.Range("C179").Value = 0
$ws.Range("A180").Value = "Republica del Chad"
$ws.Range("A181").Value = "Republica de Africa Central"
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 1
$ws.Range("A182").Value = "Nepal"
$ws.Range("B182").Value = 1
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 1
$ws.Range("E182").Value = 0
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0
